$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 286882.75
$ws.Range("I6").Value = 833441.8
$ws.Range("J6").Value = 68259.13
$ws.Range("K6").Value = 2500325.4
$ws.Range("L6").Value = 204777.39
$ws.Range("M6").Value = -2500213.4
$ws.Range("N6").Value = -205001.39
# Row 12
$ws.Range("H12").Value = 218.4
$ws.Range("J12").Value = 200
$ws.Range("L12").Value = 200
$ws.Range("N12").Value = -540
# Row 42
$ws.Range("H42").Value = 91.75
$ws.Range("I42").Value = 88.2
$ws.Range("J42").Value = 97.666664
$ws.Range("K42").Value = 264.6
$ws.Range("L42").Value = 292.999992
$ws.Range("M42").Value = -34.60000000000002
$ws.Range("N42").Value = -752.999992
# Row 62
$ws.Range("H62").Value = 200008000
$ws.Range("I62").Value = 333340000
$ws.Range("K62").Value = 333340000
$ws.Range("M62").Value = -333339376
# Row 65
$ws.Range("H65").Value = 200008000
$ws.Range("I65").Value = 333340000
$ws.Range("K65").Value = 1666700000
$ws.Range("M65").Value = -1666696880
# Row 70
$ws.Range("H70").Value = 28573356
$ws.Range("I70").Value = 200000000
$ws.Range("J70").Value = 2249.3333
$ws.Range("K70").Value = 600000000
$ws.Range("L70").Value = 6747.999899999999
$ws.Range("M70").Value = -599999730
$ws.Range("N70").Value = -7287.999899999999
# Row 73
$ws.Range("H73").Value = 28573356
$ws.Range("I73").Value = 200000000
$ws.Range("J73").Value = 2249.3333
$ws.Range("K73").Value = 600000000
$ws.Range("L73").Value = 6747.999899999999
$ws.Range("M73").Value = -599999064
$ws.Range("N73").Value = -8619.999899999999
# Row 92
$ws.Range("H92").Value = 1537.6666
$ws.Range("I92").Value = 1486.5454
$ws.Range("J92").Value = 2100
$ws.Range("K92").Value = 1486.5454
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -238.5454
$ws.Range("N92").Value = -4596
# Row 98
$ws.Range("H98").Value = 24002.4
$ws.Range("I98").Value = 25524.305
$ws.Range("J98").Value = 6500.5
$ws.Range("K98").Value = 25524.305
$ws.Range("L98").Value = 6500.5
$ws.Range("M98").Value = -24026.305
$ws.Range("N98").Value = -9496.5
# Row 113
$ws.Range("H113").Value = 16384
$ws.Range("I113").Value = 22499
$ws.Range("J113").Value = 13326.5
$ws.Range("K113").Value = 22499
$ws.Range("L113").Value = 13326.5
$ws.Range("M113").Value = -19245
$ws.Range("N113").Value = -19834.5
# Row 122
$ws.Range("H122").Value = 24002.4
$ws.Range("I122").Value = 25524.305
$ws.Range("J122").Value = 6500.5
$ws.Range("K122").Value = 76572.91500000001
$ws.Range("L122").Value = 19501.5
$ws.Range("M122").Value = -74122.91500000001
$ws.Range("N122").Value = -24401.5
# Row 135
$ws.Range("H135").Value = 16732.166
$ws.Range("J135").Value = 12592.333
$ws.Range("L135").Value = 113330.997
$ws.Range("N135").Value = -118400.997
# Row 138
$ws.Range("H138").Value = 283565.03
$ws.Range("I138").Value = 471528.03
$ws.Range("K138").Value = 1414584.09
$ws.Range("M138").Value = -1409444.09
# Row 141
$ws.Range("H141").Value = 5950.5386
$ws.Range("I141").Value = 5422.5654
$ws.Range("J141").Value = 9998.333000000001
$ws.Range("K141").Value = 16267.6962
$ws.Range("L141").Value = 29994.999
$ws.Range("M141").Value = -11087.6962
$ws.Range("N141").Value = -40354.999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7942.164
$ws.Range("I32").Value = 7942.164
$ws.Range("K32").Value = 7942.164
$ws.Range("M32").Value = -7655.164
# Row 74
$ws.Range("H74").Value = 3141.5325
$ws.Range("I74").Value = 7065.35
$ws.Range("K74").Value = 7065.35
$ws.Range("M74").Value = -6191.35
# Row 77
$ws.Range("H77").Value = 3141.5325
$ws.Range("I77").Value = 7065.35
$ws.Range("K77").Value = 35326.75
$ws.Range("M77").Value = -30958.75
# Row 122
$ws.Range("H122").Value = 1205749.8
$ws.Range("I122").Value = 5512.8423
$ws.Range("J122").Value = 5006500
$ws.Range("K122").Value = 16538.5269
$ws.Range("L122").Value = 15019500
$ws.Range("M122").Value = -14088.5269
$ws.Range("N122").Value = -15024400
# Row 132
$ws.Range("H132").Value = 3001.5454
$ws.Range("I132").Value = 1431.2858
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 4293.857400000001
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -1763.857400000001
$ws.Range("N132").Value = -22308.5
# Row 135
$ws.Range("H135").Value = 199985.53
$ws.Range("J135").Value = 199985.53
$ws.Range("L135").Value = 199985.53
$ws.Range("N135").Value = -210125.53
# Row 138
$ws.Range("H138").Value = 68991.5
$ws.Range("J138").Value = 68991.5
$ws.Range("L138").Value = 68991.5
$ws.Range("N138").Value = -79271.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 3321.394
$ws.Range("I94").Value = 2237.3704
$ws.Range("K94").Value = 2237.3704
$ws.Range("M94").Value = -1786.3704
# Row 107
$ws.Range("H107").Value = 2511
$ws.Range("I107").Value = 2511
$ws.Range("K107").Value = 2511
$ws.Range("M107").Value = -591
# Row 117
$ws.Range("H117").Value = 19666.334
$ws.Range("J117").Value = 19666.334
$ws.Range("L117").Value = 19666.334
$ws.Range("N117").Value = -28844.334
# Row 134
$ws.Range("H134").Value = 2087.5938
$ws.Range("I134").Value = 1544.931
$ws.Range("J134").Value = 7333.3335
$ws.Range("K134").Value = 4634.793
$ws.Range("L134").Value = 22000.0005
$ws.Range("M134").Value = -2099.793
$ws.Range("N134").Value = -27070.0005

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3186.1794
$ws.Range("I31").Value = 2105.4546
$ws.Range("J31").Value = 4584.7646
$ws.Range("K31").Value = 2105.4546
$ws.Range("L31").Value = 4584.7646
$ws.Range("M31").Value = -1810.4546
$ws.Range("N31").Value = -5174.7646
# Row 34
$ws.Range("H34").Value = 3186.1794
$ws.Range("I34").Value = 2105.4546
$ws.Range("J34").Value = 4584.7646
$ws.Range("K34").Value = 2105.4546
$ws.Range("L34").Value = 4584.7646
$ws.Range("M34").Value = -1903.4546
$ws.Range("N34").Value = -4988.7646
# Row 70
$ws.Range("H70").Value = 42867.5
$ws.Range("J70").Value = 42867.5
$ws.Range("L70").Value = 42867.5
$ws.Range("N70").Value = -43497.5
# Row 73
$ws.Range("H73").Value = 42867.5
$ws.Range("J73").Value = 42867.5
$ws.Range("L73").Value = 42867.5
$ws.Range("N73").Value = -45051.5
# Row 93
$ws.Range("H93").Value = 7504.2
$ws.Range("I93").Value = 7504.2
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 7504.2
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -5632.2
$ws.Range("N93").ClearContents()
# Row 99
$ws.Range("H99").Value = 197021.16
$ws.Range("I99").Value = 360710.78
$ws.Range("K99").Value = 360710.78
$ws.Range("M99").Value = -359212.78
# Row 122
$ws.Range("H122").Value = 60777.5
$ws.Range("I122").Value = 60777.5
$ws.Range("K122").Value = 182332.5
$ws.Range("M122").Value = -179882.5
# Row 126
$ws.Range("H126").Value = 197021.16
$ws.Range("I126").Value = 360710.78
$ws.Range("K126").Value = 1082132.34
$ws.Range("M126").Value = -1079662.34
# Row 132
$ws.Range("H132").Value = 10885.546
$ws.Range("I132").Value = 12584.223
$ws.Range("K132").Value = 37752.669
$ws.Range("M132").Value = -35222.669
# Row 134
$ws.Range("H134").Value = 2237.2307
$ws.Range("I134").Value = 1749.6857
$ws.Range("K134").Value = 5249.0571
$ws.Range("M134").Value = -2714.0571
# Row 141
$ws.Range("H141").Value = 313744.34
$ws.Range("J141").Value = 349961.94
$ws.Range("L141").Value = 349961.94
$ws.Range("N141").Value = -360321.94

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 121
$ws.Range("I7").Value = 115.5
$ws.Range("J7").Value = 140.8
$ws.Range("K7").Value = 346.5
$ws.Range("L7").Value = 422.4
$ws.Range("M7").Value = -234.5
$ws.Range("N7").Value = -646.4000000000001
# Row 23
$ws.Range("H23").Value = 9804694
$ws.Range("I23").Value = 659.75
$ws.Range("J23").Value = 18519390
$ws.Range("K23").Value = 1979.25
$ws.Range("L23").Value = 55558170
$ws.Range("M23").Value = -1744.25
$ws.Range("N23").Value = -55558640
# Row 33
$ws.Range("H33").Value = 229.8
$ws.Range("J33").Value = 354.66666
$ws.Range("L33").Value = 2127.99996
$ws.Range("N33").Value = -2693.99996
# Row 40
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 97
$ws.Range("H97").Value = 19655.646
$ws.Range("I97").Value = 24112.77
$ws.Range("J97").Value = 5170
$ws.Range("K97").Value = 72338.31
$ws.Range("L97").Value = 15510
$ws.Range("M97").Value = -71842.31
$ws.Range("N97").Value = -16502

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 8892.951999999999
$ws.Range("I102").Value = 9597.333000000001
$ws.Range("K102").Value = 9597.333000000001
$ws.Range("M102").Value = -7975.333000000001
# Row 122
$ws.Range("H122").Value = 7789.6855
$ws.Range("I122").Value = 6084.05
$ws.Range("J122").Value = 10063.866
$ws.Range("K122").Value = 18252.15
$ws.Range("L122").Value = 30191.598
$ws.Range("M122").Value = -15802.15
$ws.Range("N122").Value = -35091.598
# Row 126
$ws.Range("H126").Value = 11123.134
$ws.Range("J126").Value = 4399.6
$ws.Range("L126").Value = 13198.8
$ws.Range("N126").Value = -18138.8
# Row 132
$ws.Range("H132").Value = 3604.7942
$ws.Range("I132").Value = 3809.5
$ws.Range("J132").Value = 3374.5
$ws.Range("K132").Value = 11428.5
$ws.Range("L132").Value = 10123.5
$ws.Range("M132").Value = -8898.5
$ws.Range("N132").Value = -15183.5
# Row 134
$ws.Range("H134").Value = 90162.5
$ws.Range("J134").Value = 90162.5
$ws.Range("L134").Value = 270487.5
$ws.Range("N134").Value = -275557.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 748.7619
$ws.Range("J22").Value = 839.5714
$ws.Range("L22").Value = 839.5714
$ws.Range("N22").Value = -1429.5714
# Row 27
$ws.Range("H27").Value = 748.7619
$ws.Range("J27").Value = 839.5714
$ws.Range("L27").Value = 839.5714
$ws.Range("N27").Value = -1053.5714
# Row 40
$ws.Range("H40").Value = 28758.61
$ws.Range("I40").Value = 49345.2
$ws.Range("K40").Value = 49345.2
$ws.Range("M40").Value = -49209.2
# Row 46
$ws.Range("H46").Value = 2313.2666
$ws.Range("I46").Value = 812.5
$ws.Range("K46").Value = 812.5
$ws.Range("M46").Value = -624.5
# Row 61
$ws.Range("H61").Value = 33994.75
$ws.Range("I61").Value = 2112.889
$ws.Range("J61").Value = 74985.71000000001
$ws.Range("K61").Value = 2112.889
$ws.Range("L61").Value = 74985.71000000001
$ws.Range("M61").Value = -1910.889
$ws.Range("N61").Value = -75389.71000000001
# Row 64
$ws.Range("H64").Value = 35800
$ws.Range("J64").Value = 35800
$ws.Range("L64").Value = 35800
$ws.Range("N64").Value = -36250
# Row 67
$ws.Range("H67").Value = 35800
$ws.Range("J67").Value = 35800
$ws.Range("L67").Value = 35800
$ws.Range("N67").Value = -37360
# Row 68
$ws.Range("H68").Value = 5277.75
$ws.Range("J68").Value = 6542.769
$ws.Range("L68").Value = 6542.769
$ws.Range("N68").Value = -8040.769
# Row 71
$ws.Range("H71").Value = 5277.75
$ws.Range("J71").Value = 6542.769
$ws.Range("L71").Value = 32713.845
$ws.Range("N71").Value = -40201.845
# Row 93
$ws.Range("H93").Value = 4557.75
$ws.Range("I93").Value = 4564.278
$ws.Range("J93").Value = 4499
$ws.Range("K93").Value = 4564.278
$ws.Range("L93").Value = 4499
$ws.Range("M93").Value = -3316.278
$ws.Range("N93").Value = -6995
# Row 96
$ws.Range("H96").Value = 22500
$ws.Range("J96").Value = 22500
$ws.Range("L96").Value = 22500
$ws.Range("N96").Value = -27992
# Row 100
$ws.Range("H100").Value = 6324.5
$ws.Range("I100").Value = 4861.875
$ws.Range("J100").Value = 9249.75
$ws.Range("K100").Value = 4861.875
$ws.Range("L100").Value = 9249.75
$ws.Range("M100").Value = -4320.875
$ws.Range("N100").Value = -10331.75
# Row 113
$ws.Range("H113").Value = 33994.75
$ws.Range("I113").Value = 2112.889
$ws.Range("J113").Value = 74985.71000000001
$ws.Range("K113").Value = 2112.889
$ws.Range("L113").Value = 74985.71000000001
$ws.Range("M113").Value = 57.11099999999988
$ws.Range("N113").Value = -79325.71000000001
# Row 122
$ws.Range("H122").Value = 7841.1665
$ws.Range("I122").Value = 9306.286
$ws.Range("J122").Value = 6908.8184
$ws.Range("K122").Value = 27918.858
$ws.Range("L122").Value = 20726.4552
$ws.Range("M122").Value = -25468.858
$ws.Range("N122").Value = -25626.4552
# Row 132
$ws.Range("H132").Value = 1153544.1
$ws.Range("I132").Value = 2137475
$ws.Range("J132").Value = 5624.6665
$ws.Range("K132").Value = 6412425
$ws.Range("L132").Value = 16873.9995
$ws.Range("M132").Value = -6409895
$ws.Range("N132").Value = -21933.9995
# Row 136
$ws.Range("H136").Value = 5055.7095
$ws.Range("I136").Value = 3414.9333
$ws.Range("J136").Value = 6593.9375
$ws.Range("K136").Value = 10244.7999
$ws.Range("L136").Value = 19781.8125
$ws.Range("M136").Value = -7694.7999
$ws.Range("N136").Value = -24881.8125

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 297833.5
$ws.Range("I62").Value = 436874.88
$ws.Range("J62").Value = 19750.75
$ws.Range("K62").Value = 436874.88
$ws.Range("L62").Value = 19750.75
$ws.Range("M62").Value = -436250.88
$ws.Range("N62").Value = -20998.75
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 65
$ws.Range("H65").Value = 297833.5
$ws.Range("I65").Value = 436874.88
$ws.Range("J65").Value = 19750.75
$ws.Range("K65").Value = 2184374.4
$ws.Range("L65").Value = 98753.75
$ws.Range("M65").Value = -2181254.4
$ws.Range("N65").Value = -104993.75
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 122
$ws.Range("H122").Value = 4904.375
$ws.Range("I122").Value = 2917.8948
$ws.Range("J122").Value = 7807.6924
$ws.Range("K122").Value = 8753.6844
$ws.Range("L122").Value = 23423.0772
$ws.Range("M122").Value = -6303.6844
$ws.Range("N122").Value = -28323.0772
# Row 132
$ws.Range("H132").Value = 4318.098
$ws.Range("I132").Value = 4371.114
$ws.Range("J132").Value = 4180.8823
$ws.Range("K132").Value = 13113.342
$ws.Range("L132").Value = 12542.6469
$ws.Range("M132").Value = -10583.342
$ws.Range("N132").Value = -17602.6469
# Row 135
$ws.Range("H135").Value = 112777.5
$ws.Range("J135").Value = 155555
$ws.Range("L135").Value = 155555
$ws.Range("N135").Value = -165695
# Row 136
$ws.Range("H136").Value = 361321.78
$ws.Range("I136").Value = 419314.78
$ws.Range("K136").Value = 1257944.34
$ws.Range("M136").Value = -1255394.34

